# Update "provincias_spain" COVID dashboard sheet:
#  - refresh the "last updated" timestamp
#  - re-sort / update a handful of provinces (Navarra, Sevilla, La Rioja,
#    Alacant/Alicante, Albacete, Araba/Alava, Castello/Castellon, Guadalajara)
#  - update the numeric totals for several rows (Madrid, Valencia/Valencia,
#    Navarra, Sevilla, La Rioja, Alacant/Alicante, Albacete, Araba/Alava,
#    Cantabria, Castello/Castellon, Guadalajara)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp update (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 12:22"

# --- Province name updates (city re-ordering) ---
$ws.Range("A11").Value = "Navarra"
$ws.Range("A12").Value = "Sevilla"
$ws.Range("A15").Value = "La Rioja"
$ws.Range("A16").Value = "Alacant/Alicante"
$ws.Range("A17").Value = "Albacete"
$ws.Range("A18").Value = "Araba/Alava"
$ws.Range("A38").Value = "Castello/Castellon"
$ws.Range("A39").Value = "Guadalajara"

# --- Row 4: Madrid ---
$ws.Range("B4").Value = 38723
$ws.Range("C4").Value = 17322
$ws.Range("D4").Value = 16265
$ws.Range("E4").Value = 5136

# --- Row 9: Valencia/Valencia ---
$ws.Range("B9").Value = 3788
$ws.Range("C9").Value = 561
$ws.Range("D9").Value = 2929
$ws.Range("E9").Value = 298

# --- Row 11: Navarra ---
$ws.Range("B11").Value = 3231
$ws.Range("C11").Value = 380
$ws.Range("D11").Value = 2664
$ws.Range("E11").Value = 187

# --- Row 12: Sevilla ---
$ws.Range("B12").Value = 3137
$ws.Range("C12").Value = 82
$ws.Range("D12").Value = 2962
$ws.Range("E12").Value = 93

# --- Row 15: La Rioja ---
$ws.Range("B15").Value = 2719
$ws.Range("C15").Value = 964
$ws.Range("D15").Value = 1614
$ws.Range("E15").Value = 141

# --- Row 16: Alacant/Alicante ---
$ws.Range("B16").Value = 2673
$ws.Range("C16").Value = 284
$ws.Range("D16").Value = 2120
$ws.Range("E16").Value = 269

# --- Row 17: Albacete ---
$ws.Range("B17").Value = 2653
$ws.Range("C17").Value = 1149
$ws.Range("D17").Value = 7827
$ws.Range("E17").Value = 208

# --- Row 18: Araba/Alava ---
$ws.Range("B18").Value = 2639
$ws.Range("C18").Value = 3405
$ws.Range("D18").Value = 4708
$ws.Range("E18").Value = 194

# --- Row 29: Cantabria ---
$ws.Range("B29").Value = 1483
$ws.Range("C29").Value = 118
$ws.Range("D29").Value = 1288
$ws.Range("E29").Value = 77

# --- Row 38: Castello/Castellon ---
$ws.Range("B38").Value = 869
$ws.Range("C38").Value = 85
$ws.Range("D38").Value = 714
$ws.Range("E38").Value = 70

# --- Row 39: Guadalajara ---
$ws.Range("B39").Value = 858
$ws.Range("C39").Value = 1149
$ws.Range("D39").Value = 7827
$ws.Range("E39").Value = 121
